# Add 2 more DM levels to show variation in output.
# Sheet "Slurry": update C2:C4 (and C5, which becomes part of the first
# DM-level block) from 5.9 -> 5.1, then append two more copies of the
# original 4-row block (DM=5.9 and DM=6.9) as rows 6-9 and rows 10-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# --- Update existing rows 2-5: man.dm (column C) changes from 5.9 to 5.1 ---
$ws.Cells.Item(2, 3).Value = 5.1
$ws.Cells.Item(3, 3).Value = 5.1
$ws.Cells.Item(4, 3).Value = 5.1
$ws.Cells.Item(5, 3).Value = 5.1

# --- New block 1 (rows 6-9): DM = 5.9, same layout as the original block ---
$ws.Cells.Item(6, 1).Value = $ws.Cells.Item(2, 1).Value2
$ws.Cells.Item(6, 2).Value = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(6, 3).Value = 5.9
$ws.Cells.Item(6, 4).Value = 7.9

$ws.Cells.Item(7, 1).Value = $ws.Cells.Item(3, 1).Value2
$ws.Cells.Item(7, 2).Value = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(7, 3).Value = 5.9
$ws.Cells.Item(7, 4).Formula = "=7.9-1.38"

$ws.Cells.Item(8, 1).Value = $ws.Cells.Item(4, 1).Value2
$ws.Cells.Item(8, 2).Value = $ws.Cells.Item(4, 2).Value2
$ws.Cells.Item(8, 3).Value = 5.9
$ws.Cells.Item(8, 4).Formula = "=7.9-0.8187"

$ws.Cells.Item(9, 1).Value = $ws.Cells.Item(5, 1).Value2
$ws.Cells.Item(9, 2).Value = $ws.Cells.Item(5, 2).Value2
$ws.Cells.Item(9, 3).Value = 5.9
$ws.Cells.Item(9, 4).Formula = "=7.9-1.11"

# --- New block 2 (rows 10-13): DM = 6.9, same layout as the original block ---
$ws.Cells.Item(10, 1).Value = $ws.Cells.Item(2, 1).Value2
$ws.Cells.Item(10, 2).Value = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(10, 3).Value = 6.9
$ws.Cells.Item(10, 4).Value = 7.9

$ws.Cells.Item(11, 1).Value = $ws.Cells.Item(3, 1).Value2
$ws.Cells.Item(11, 2).Value = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(11, 3).Value = 6.9
$ws.Cells.Item(11, 4).Formula = "=7.9-1.38"

$ws.Cells.Item(12, 1).Value = $ws.Cells.Item(4, 1).Value2
$ws.Cells.Item(12, 2).Value = $ws.Cells.Item(4, 2).Value2
$ws.Cells.Item(12, 3).Value = 6.9
$ws.Cells.Item(12, 4).Formula = "=7.9-0.8187"

$ws.Cells.Item(13, 1).Value = $ws.Cells.Item(5, 1).Value2
$ws.Cells.Item(13, 2).Value = $ws.Cells.Item(5, 2).Value2
$ws.Cells.Item(13, 3).Value = 6.9
$ws.Cells.Item(13, 4).Formula = "=7.9-1.11"

# --- Copy the style/number-format of the template rows onto the new rows ---
$ws.Range("A2:D2").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A10:D10").PasteSpecial(-4122)

$ws.Range("A3:D5").Copy()
$ws.Range("A7:D9").PasteSpecial(-4122)
$ws.Range("A11:D13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Move the selection to D16, matching the saved selection in the diff ---
$ws.Range("D16").Select()
